# Weekly data update: insert a new price observation for
# "Vega Modelo de Temuco - Coco" right above the existing row 25,
# pushing all subsequent rows down by one (dimension grows from
# A1:T130 to A1:T131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25 - everything currently at
# row 25..130 shifts down to 26..131.
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new weekly observation. Most
# descriptive columns repeat the values already used for this
# market/product/quality combination; only the date and the price
# columns carry new figures.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 45250
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108007
$ws.Range("J25").Value = "Coco"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 15
$ws.Range("N25").Value = 34000
$ws.Range("O25").Value = 34000
$ws.Range("P25").Value = 34000
$ws.Range("Q25").Value = "$/malla 20 unidades"
$ws.Range("R25").Value = "Perú"
$ws.Range("S25").Value = 1700
$ws.Range("T25").Value = 20
